$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows 23-26 appended below the existing table (previously A1:J22).
# Columns A (ID_Set) and C (nbPieces) must be stored as text, consistent with
# the rest of the sheet, so force a Text number format before assigning values,
# then clear the format again so no extra style is left applied to the cells.
$ws.Range("A23:A26").NumberFormat = "@"
$ws.Range("C23:C26").NumberFormat = "@"

# Row 23
$ws.Range("A23").Value = "10370"
$ws.Range("B23").Value = "L’étoile de Noël"
$ws.Range("C23").Value = "608"
$ws.Range("D23").Value = "The Botanical Collection"
$ws.Range("E23").Value = "https://www.lego.com/cdn/cs/set/assets/blt5b9064fcb12ba88b/10370_Prod.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1"
$ws.Range("F23").Value = "https://www.lego.com/fr-fr/product/10370"

# Row 24
$ws.Range("A24").Value = "43278"
$ws.Range("B24").Value = "Château d’Arendelle et palais de glace d’Elsa miniatures"
$ws.Range("C24").Value = "306"
$ws.Range("D24").Value = "Disney™"
$ws.Range("E24").Value = "https://www.lego.com/cdn/cs/set/assets/blt9b30046d62bfedb3/43278_Prod_en-gb.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1"
$ws.Range("F24").Value = "https://www.lego.com/fr-fr/product/43278"

# Row 25
$ws.Range("A25").Value = "40478"
$ws.Range("B25").Value = "Le château Disney miniature"
$ws.Range("C25").Value = "567"
$ws.Range("D25").Value = "Disney™"
$ws.Range("E25").Value = "https://www.lego.com/cdn/cs/set/assets/blt9a53be3e8553bce6/40478_Prod.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1"
$ws.Range("F25").Value = "https://www.lego.com/fr-fr/product/40478"
$ws.Range("J25").Value = "https://www.avenuedelabrique.com/lego-disney/40478-le-chateau-disney-miniature/p7596"

# Row 26
$ws.Range("A26").Value = "43260"
$ws.Range("B26").Value = "L'île de Vaiana miniature"
$ws.Range("C26").Value = "175"
$ws.Range("D26").Value = "Disney™"
$ws.Range("E26").Value = "https://www.lego.com/cdn/cs/set/assets/blte3899e6f083ec338/43260_Prod.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1"
$ws.Range("F26").Value = "https://www.lego.com/fr-fr/product/43260"

# Restore default (general) formatting for the numeric-looking text columns
# now that the values have been stored as text.
$ws.Range("A23:A26").ClearFormats()
$ws.Range("C23:C26").ClearFormats()
